# Adds a "Werkelijk resultaat" (Actual result) row to the column-mapping
# overview, right after "Verwacht resultaat" (row 15 -> new row 16),
# pushing the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 16
$ws.Rows("16:16").Insert()

# Copy the formatting of the row above ("Verwacht resultaat") onto the new
# row so it matches the rest of the table (borders/fills per column)
$ws.Range("A15:G15").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row's text
$ws.Range("A16").Value = "Werkelijk resultaat"
$ws.Range("B16").Value = "Hybrid"
$ws.Range("C16").Value = "Werkelijk resultaat"
$ws.Range("G16").Value = "Werkelijk resultaat"

# Match the author's last active selection
[void]$ws.Range("C12").Select()
